# Apply the commit: update report date headers and swap a set of
# mismatched stock-report row pairs (columns B, E, F, G) back into the
# correct order.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update "From :" (I1) and "To :" (K1) dates -------------------------
$ws.Range("I1").Value = "07-12-2025 00:00:00"
$ws.Range("K1").Value = "07-12-2025 00:00:00"

# --- Swap B/E/F/G values between each pair of rows -----------------------
$rowPairs = @(
    @(227,228),
    @(229,230),
    @(232,233),
    @(243,244),
    @(366,367),
    @(372,373),
    @(375,376),
    @(380,381),
    @(463,464),
    @(473,474),
    @(572,573)
)

$cols = @("B","E","F","G")

foreach ($pair in $rowPairs) {
    $r1 = $pair[0]
    $r2 = $pair[1]
    foreach ($col in $cols) {
        $cell1 = $ws.Range("$col$r1")
        $cell2 = $ws.Range("$col$r2")
        $tmp = $cell1.Value2
        $cell1.Value2 = $cell2.Value2
        $cell2.Value2 = $tmp
    }
}
